$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9768130779266357
$ws.Range("B1").Value = 1.743854641914368
$ws.Range("C1").Value = 5.775450706481934
$ws.Range("D1").Value = 3.43626594543457
$ws.Range("E1").Value = 0.426500529050827
